$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item(1)
$wsVersion = $wb.Worksheets.Item(2)

# --- 1. Bump the schema version (shared string index 1: "1" -> "2") ---
$wsVersion.Range("A1").Value = "2"

# --- 2. Insert two new columns (concentration_value, concentration_unit)
#        before the existing "conjugated_cat_number" column (H) ---
$wsData.Range("H1:I1").EntireColumn.Insert()

$wsData.Range("H1").Value = "concentration_value"
$wsData.Range("I1").Value = "concentration_unit"

# carry over the style used by the rest of the header row
$wsData.Range("H1:I1").Value2 = $wsData.Range("H1:I1").Value2
$headerStyle = $wsData.Range("G1").Style
$wsData.Range("H1").Style = $headerStyle
$wsData.Range("I1").Style = $headerStyle

Write-Host "done"
